# Update column widths and header labels on every "PiN" sheet of the workbook.
$wb = $excel.ActiveWorkbook

# The ColumnWidth property (character units) is persisted to the OOXML
# <col width="..."> attribute with a constant +5/6 padding baked in by the
# runtime, so subtract it here to land on the exact target width.
$widthPad = 5 / 6

# New column widths for columns I..R (9..18), expressed as the desired
# final OOXML <col width> values.
$colWidths = @{
    9  = 23   # I
    10 = 23   # J
    11 = 20   # K
    12 = 20   # L
    13 = 20   # M
    14 = 20   # N
    15 = 20   # O
    16 = 20   # P
    17 = 33   # Q
    18 = 33   # R
}

# New header labels for row 5, columns I..R
$headerText = @{
    "I5" = "% severity levels 1-2"
    "J5" = "# severity levels 1-2"
    "K5" = "% severity level 3"
    "L5" = "# severity level 3"
    "M5" = "% severity level 4"
    "N5" = "# severity level 4"
    "O5" = "% severity level 5"
    "P5" = "# severity level 5"
    "Q5" = "% Tot PiN (severity levels 3-5)"
    "R5" = "# Tot PiN (severity levels 3-5)"
}

foreach ($ws in $wb.Worksheets) {
    foreach ($colIndex in $colWidths.Keys) {
        $ws.Columns.Item($colIndex).ColumnWidth = $colWidths[$colIndex] - $widthPad
    }

    foreach ($cellRef in $headerText.Keys) {
        $ws.Range($cellRef).Value = $headerText[$cellRef]
    }
}
